# regulation-import.xlsx — add a "Mô tả" (Description) column.
#
# The sheet currently looks like:
#   A: STT | B: Tên quy định | C: Điểm trừ | D: Loại trừ điểm | E: Tiêu chí
# A new "Mô tả" column is inserted right after "Tên quy định" (i.e. becomes
# the new column C), pushing the remaining columns one place to the right:
#   A: STT | B: Tên quy định | C: Mô tả | D: Điểm trừ | E: Loại trừ điểm | F: Tiêu chí

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before the current column C ("Điểm trừ"), shifting
# C:E to D:F.
$ws.Columns("C:C").Insert() | Out-Null

# Give the new column the same width as column B, matching the look of the
# neighbouring "Tên quy định" column.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth()

# Fill in the new header and the two data rows.
$ws.Range("C1").Value = "Mô tả"
$ws.Range("C2").Value = "Mô tả 1"
$ws.Range("C3").Value = "Mô tả 2"

# Match the author's final cursor position.
$ws.Range("C4").Select() | Out-Null
